$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 430, pushing the existing
# rows 430-435 down to 432-437.
$ws.Rows.Item(430).Insert()
$ws.Rows.Item(431).Insert()

# New row 430: Limon, 1a plateado, Provincia de Melipilla, 18kg mesh
$ws.Cells.Item(430, 1).Value = 4
$ws.Cells.Item(430, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(430, 3).Value = "Los Lagos"
$ws.Cells.Item(430, 4).Value = 44628
$ws.Cells.Item(430, 5).Value = 10
$ws.Cells.Item(430, 6).Value = "Fruta"
$ws.Cells.Item(430, 7).Value = 100102
$ws.Cells.Item(430, 8).Value = "Cítricos"
$ws.Cells.Item(430, 9).Value = 100102003
$ws.Cells.Item(430, 10).Value = "Limón"
$ws.Cells.Item(430, 11).Value = "Sin especificar"
$ws.Cells.Item(430, 12).Value = "1a plateado"
$ws.Cells.Item(430, 13).Value = 1000
$ws.Cells.Item(430, 14).Value = 28000
$ws.Cells.Item(430, 15).Value = 28000
$ws.Cells.Item(430, 16).Value = 28000
$ws.Cells.Item(430, 17).Value = "$/malla 18 kilos"
$ws.Cells.Item(430, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(430, 19).Value = 1556
$ws.Cells.Item(430, 20).Value = 18

# New row 431: Limon, 2a plateado, Provincia de Melipilla, 18kg mesh
$ws.Cells.Item(431, 1).Value = 4
$ws.Cells.Item(431, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(431, 3).Value = "Los Lagos"
$ws.Cells.Item(431, 4).Value = 44628
$ws.Cells.Item(431, 5).Value = 10
$ws.Cells.Item(431, 6).Value = "Fruta"
$ws.Cells.Item(431, 7).Value = 100102
$ws.Cells.Item(431, 8).Value = "Cítricos"
$ws.Cells.Item(431, 9).Value = 100102003
$ws.Cells.Item(431, 10).Value = "Limón"
$ws.Cells.Item(431, 11).Value = "Sin especificar"
$ws.Cells.Item(431, 12).Value = "2a plateado"
$ws.Cells.Item(431, 13).Value = 500
$ws.Cells.Item(431, 14).Value = 25000
$ws.Cells.Item(431, 15).Value = 25000
$ws.Cells.Item(431, 16).Value = 25000
$ws.Cells.Item(431, 17).Value = "$/malla 18 kilos"
$ws.Cells.Item(431, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(431, 19).Value = 1389
$ws.Cells.Item(431, 20).Value = 18
